$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C ("Industry"), shifting the existing
# Mutual Fund/Status/Jan_2026/Dec_2025/Oct_2025/MoM/QoQ columns from C:I to D:J.
$ws.Columns("C").Insert()

# Header
$ws.Cells.Item(1, 3).Value = "Industry"

# Industry values for each ISIN/Stock row (rows 2-24)
$industries = @(
    "Leisure Services",
    "Finance",
    "Food Products",
    "Personal Products",
    "Automobiles",
    "Retailing",
    "Pharmaceuticals & Biotechnology",
    "Consumer Durables",
    "Agricultural Food & other Products",
    "Personal Products",
    "Consumer Durables",
    "Finance",
    "Banks",
    "Beverages",
    "Consumer Durables",
    "Beverages",
    "Retailing",
    "IT - Software",
    "Realty",
    "Power",
    "Agricultural Food & other Products",
    "Retailing",
    "N.A."
)

for ($i = 0; $i -lt $industries.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $industries[$i]
}
